$d = $word.ActiveDocument

# Step 1: Replace "hello" with "H" (keeps the original run's formatting).
$d.Content.Find.Execute("hello", $true, $false, $false, $false, $false, $true, 1, $false, "H", 2)

# Step 2: Append "ello" and then ", how are you?" right after "H", so the
# paragraph text becomes "Hello, how are you?". Because the appended text
# shares identical run formatting with what precedes it, Word merges it
# into a single run at this point - that's fine, we'll split it below.
$para = $d.Paragraphs.Item(1).Range
$insertPoint = $d.Range($para.End - 1, $para.End - 1)
$insertPoint.InsertAfter("ello")

$para = $d.Paragraphs.Item(1).Range
$insertPoint2 = $d.Range($para.End - 1, $para.End - 1)
$insertPoint2.InsertAfter(", how are you?")

# Step 3: Split "Hello, how are you?" into three separate runs: "H",
# "ello", and ", how are you?" - matching the target structure. We force
# each piece into its own <w:r> by toggling a direct-formatting property
# (Bold) on and back off; this keeps the runs from being re-merged on
# save even though their final formatting ends up identical, while not
# changing how the text actually looks. Go from the end of the paragraph
# towards the start so earlier edits don't shift the offsets we still
# need to use.
$restRange = $d.Range(5, 19)
$restRange.Bold = 1
$restRange.Bold = 0

$elloRange = $d.Range(1, 5)
$elloRange.Bold = 1
$elloRange.Bold = 0
